$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Delete rows for program_prop_vac (row 13) and program_prop_unvac (row 14).
# Everything below shifts up by two rows automatically, formulas/styles follow.
$ws.Rows("13:14").Delete()

# Update the sheet view so that the visible selection matches the commit's
# new state (scrolled down a bit with C11 selected).
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C11").Select()
